$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64, shifting existing rows 64-172 down to 65-173.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record.
$ws.Range("A64").Value = 10
$ws.Range("B64").Value = "Vega Modelo de Temuco"
$ws.Range("C64").Value = "La Araucanía"
$ws.Range("D64").Value = 44536
$ws.Range("E64").Value = 9
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100102
$ws.Range("H64").Value = "Cítricos"
$ws.Range("I64").Value = 100102006
$ws.Range("J64").Value = "Pomelo"
$ws.Range("K64").Value = "Start Ruby"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 55
$ws.Range("N64").Value = 12000
$ws.Range("O64").Value = 12000
$ws.Range("P64").Value = 12000
$ws.Range("Q64").Value = "$/bandeja 15 kilos granel"
$ws.Range("R64").Value = "Región de O'Higgins"
$ws.Range("S64").Value = 800
$ws.Range("T64").Value = 15
